# Commit message: "remove lf from format"
#
# The "format" column of the code_roscop sheet contains C-style printf
# format strings (e.g. "%6.3f", "%9.0f", ...). A handful of them were
# written with the stray "long" length modifier "l" before the final
# conversion character (e.g. "%6.3lf" instead of "%6.3f"). This edit
# removes that extraneous "l" wherever it appears in a format string.
#
# We scan every used cell on the sheet and, whenever its text looks like
# a printf style format specifier containing "lf" (e.g. %6.3lf, %7.2lf,
# %8.4lf, ...), rewrite it to the equivalent specifier without the "l".
# Cells that don't match this pattern (dates, labels, numbers, other
# format codes such as %6.6d or %3d, etc.) are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $old = $cell.Value2

        if ($old -ne $null -and $old -is [string] -and $old -match "^%[0-9+.]*lf$") {
            $new = $old -replace "lf$", "f"
            $cell.Value = $new
        }
    }
}
